$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 ("Human Resources Associate") entirely; this shifts rows 4-5 up to 3-4
$ws.Rows.Item(3).Delete()
